# Append the latest COVID-19 data rows (2020-05-25 .. 2020-06-03) to the
# "Covid-19 podatki" sheet, growing the data table Tabela1 from A1:J75 to
# A1:J85, matching the upstream bot's daily data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows: Date(serial), Tested(all), Tested(daily), Positive(all),
# Positive(daily), Hospitalized, Intensive care, Discharged, Deaths(all), Deaths(daily)
$newRows = @(
    @(43976, 75770, 754, 1469,   0, 9, 2, 6, 108, 1),
    @(43977, 76579, 809, 1471,   2, 8, 2, 2, 108, 0),
    @(43978, 77210, 631, 1473,   2, 7, 2, 1, 108, 0),
    @(43979, 77916, 706, 1473,   0, 7, 2, 0, 108, 0),
    @(43980, 78529, 613, 1473,   0, 7, 2, 0, 108, 0),
    @(43981, 78793, 264, 1473,   0, 6, 2, 1, 108, 0),
    @(43982, 79039, 246, 1473,   0, 5, 1, 0, 109, 1),
    @(43983, 79698, 659, 1475,   2, 5, 1, 0, 109, 0),
    @(43984, 80505, 807, 1477,   2, 5, 0, 0, 109, 0),
    @(43985, 81333, 828, 1477,   0, 5, 0, 0, 109, 0)
)

$firstNewRow = 76
$lastExistingRow = 75
$lastNewRow = $firstNewRow + $newRows.Count - 1

# Carry the last existing row's formatting down across the new rows before
# writing values, so the new rows inherit the same look (number formats,
# alignment, fonts, banding) as the rest of the table.
$ws.Range("A$lastExistingRow`:J$lastExistingRow").Copy($ws.Range("A$firstNewRow`:J$lastNewRow"))

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $firstNewRow + $i
    $vals = $newRows[$i]
    for ($c = 1; $c -le 10; $c++) {
        $ws.Cells.Item($r, $c).Value2 = $vals[$c - 1]
    }
}

# Row 81's "Tested (all)" cell (B81) carries the same plain-number format as
# the rest of that row (C81:J81) rather than the thousands-separated format
# normally used in column B - mirror that one-off formatting quirk from the
# source data refresh.
$ws.Cells.Item(81, 3).Copy($ws.Cells.Item(81, 2))
$ws.Cells.Item(81, 2).Value2 = 78793

# Grow the structured table (Tabela1) to include the newly added rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:J$lastNewRow"))

# Reflect the updated view/selection like the source workbook does after
# the refresh (selection lands on the freshly appended last row).
$ws.Range("A$lastNewRow`:J$lastNewRow").Select()
